$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '55.327.11'
$ws.Range("E2").Value = '  -1.49%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.355.67'
$ws.Range("E3").Value = '  -4.88%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '477.13'
$ws.Range("E5").Value = '  -2.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.06'
$ws.Range("E6").Value = '  +0.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.630'
$ws.Range("E7").Value = '  +23.83%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.998'
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.363.08'
$ws.Range("E9").Value = '  -5.06%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0969'
$ws.Range("E10").Value = '  -0.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.46'
$ws.Range("E11").Value = '  -5.96%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.327'
$ws.Range("E12").Value = '  -1.78%  '
$ws.Range("E13").Value = '  +0.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.765.59'
$ws.Range("E14").Value = '  -5.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '55.276.61'
$ws.Range("E15").Value = '  -1.55%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.09'
$ws.Range("E16").Value = '  -4.84%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000130'
$ws.Range("E17").Value = '  -4.39%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.358.01'
$ws.Range("E18").Value = '  -5.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.56'
$ws.Range("E19").Value = '  +0.86%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '316.29'
$ws.Range("E20").Value = '  -0.68%  '
$ws.Range("E21").Value = '  -4.70%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.64'
$ws.Range("E23").Value = '  -2.89%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '56.85'
$ws.Range("E24").Value = '  -2.80%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.396'
$ws.Range("E26").Value = '  -3.89%  '
$ws.Range("E27").Value = '  -6.15%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.447.67'
$ws.Range("E28").Value = '  -5.32%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.12'
$ws.Range("E29").Value = '  -6.55%  '
$ws.Range("E30").Value = '  +0.14%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0' + [string]([char]0x2083) + '0749'
$ws.Range("E31").Value = '  -5.10%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.19'
$ws.Range("E32").Value = '  -0.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '145.89'
$ws.Range("E33").Value = '  -2.04%  '
$ws.Range("E34").Value = '  -1.97%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.11'
$ws.Range("E35").Value = '  -1.59%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.10'
$ws.Range("E36").Value = '  -3.37%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.58'
$ws.Range("E37").Value = '  -4.23%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.815'
$ws.Range("E38").Value = '  -5.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.103'
$ws.Range("E39").Value = '  +11.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '33.69'
$ws.Range("E40").Value = '  -1.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  +0.24%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.34'
$ws.Range("E42").Value = '  +0.85%  '
$ws.Range("E43").Value = '  -3.22%  '
$ws.Range("E44").Value = '  -4.67%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0519'
$ws.Range("E45").Value = '  -6.25%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.16'
$ws.Range("E46").Value = '  -0.16%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '253.51'
$ws.Range("E47").Value = '  -1.89%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0221'
$ws.Range("E48").Value = '  -3.19%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.41'
$ws.Range("E49").Value = '  -6.95%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.79'
$ws.Range("E50").Value = '  -4.24%  '
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.800.59'
$ws.Range("E51").Value = '  -4.35%  '
